$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates derived from the crypto-price refresh diff.
# Each entry is Cell, NewValue. Cells in $textCells hold numeric-looking
# strings (prices such as "0.999") that must stay text, matching the
# original inlineStr cells instead of being auto-converted to numbers.
$textCells = @("D5","D6","D7","D10","D14","D17","D21","D25","D26","D28","D30","D31","D32","D33","D34","D35","D38","D39","D42","D44","D47","D48","D50")

$updates = [ordered]@{
    "D2" = "66.978.64"
    "E2" = "  -0.44%  "
    "D3" = "3.457.60"
    "E3" = "  -1.43%  "
    "E4" = "  -0.05%  "
    "D5" = "592.72"
    "E5" = "  -0.90%  "
    "D6" = "179.70"
    "E6" = "  +2.25%  "
    "D7" = "0.608"
    "E7" = "  +3.24%  "
    "E8" = "  -0.04%  "
    "D9" = "3.454.74"
    "E9" = "  -1.53%  "
    "D10" = "0.140"
    "E10" = "  +5.72%  "
    "E11" = "  -3.08%  "
    "E12" = "  -0.27%  "
    "D13" = "4.056.21"
    "E13" = "  -1.45%  "
    "D14" = "31.68"
    "E14" = "  +3.21%  "
    "E15" = "  -0.55%  "
    "D16" = "66.979.13"
    "E16" = "  -0.48%  "
    "D17" = "0.0000176"
    "E17" = "  -1.54%  "
    "D18" = "3.456.99"
    "E18" = "  -1.53%  "
    "E19" = "  -1.42%  "
    "E20" = "  -2.45%  "
    "D21" = "389.34"
    "E21" = "  -1.33%  "
    "E22" = "  -1.22%  "
    "E23" = "  +0.07%  "
    "E24" = "  +1.27%  "
    "B25" = "Litecoin"
    "C25" = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
    "D25" = "72.02"
    "E25" = "  -2.05%  "
    "B26" = "Polygon"
    "C26" = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
    "D26" = "0.535"
    "E26" = "  -0.43%  "
    "E27" = "  -0.87%  "
    "D28" = "10.33"
    "E28" = "  +1.31%  "
    "E29" = "  -3.58%  "
    "D30" = "0.999"
    "E30" = "  +0.33%  "
    "D31" = "6.13"
    "E31" = "  -0.33%  "
    "B32" = "PancakeSwap"
    "C32" = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
    "D32" = "2.05"
    "E32" = "  -0.49%  "
    "B33" = "Fetch.AI"
    "C33" = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
    "D33" = "1.40"
    "E33" = "  -2.23%  "
    "D34" = "23.38"
    "E34" = "  -1.29%  "
    "D35" = "7.31"
    "E35" = "  -1.01%  "
    "E36" = "  -0.06%  "
    "E37" = "  -3.75%  "
    "D38" = "163.20"
    "E38" = "  -0.15%  "
    "D39" = "0.874"
    "E39" = "  -0.83%  "
    "E40" = "  +9.56%  "
    "E41" = "  -3.06%  "
    "D42" = "6.78"
    "E42" = "  -3.38%  "
    "E43" = "  -0.66%  "
    "D44" = "26.11"
    "E44" = "  -0.28%  "
    "E45" = "  -1.90%  "
    "D46" = "2.736.95"
    "E46" = "  -2.38%  "
    "D47" = "26.21"
    "E47" = "  -5.39%  "
    "D48" = "41.26"
    "E48" = "  -2.85%  "
    "E49" = "  -1.95%  "
    "D50" = "327.41"
    "E50" = "  -4.50%  "
    "E51" = "  -4.02%  "
}

foreach ($cellRef in $updates.Keys) {
    $value = $updates[$cellRef]
    $range = $ws.Range($cellRef)
    if ($textCells -contains $cellRef) {
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.Style = "Normal"
    } else {
        $range.Value = $value
    }
}

Write-Output "Applied $($updates.Count) cell updates"
